$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.039.76"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.944.57"
$ws.Range("E3").Value = "  +1.84%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "378.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.542"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0839"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.410.52"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "2.943.28"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.961"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.92%  "
$ws.Range("D18").Value = "51.096.93"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +22.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.169"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.113"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.17%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -2.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0447"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.55%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "122.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.05"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.285"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +21.53%  "
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.74%  "
$ws.Range("D48").Value = "2.030.98"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0344"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.70%  "
$ws.Range("E51").Value = "  +2.30%  "
